# Auto-generated cell updates applying the TPM re-run diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = [double]"29.663058"
$ws.Range("H2").Value = [double]"88.98917399999999"
$ws.Range("I2").Value = [double]"0.03283316554829836"
$ws.Range("J2").Value = [double]"0.03283316554829836"
$ws.Range("O2").Value = [double]"0.001498364820294181"
$ws.Range("P2").Value = [double]"0.001498364820294181"
$ws.Range("Q2").Value = [double]"0.3614641371019999"
$ws.Range("R2").Value = [double]"3.253177233918"
$ws.Range("S2").Value = [double]"4.919606019646516E-05"
$ws.Range("T2").Value = [double]"4.919606019646515E-05"
$ws.Range("G3").Value = [double]"29.663058"
$ws.Range("H3").Value = [double]"88.98917399999999"
$ws.Range("I3").Value = [double]"0.03283316554829836"
$ws.Range("J3").Value = [double]"0.03283316554829836"
$ws.Range("M3").Value = [double]"0.08128566666666666"
$ws.Range("N3").Value = [double]"0.243857"
$ws.Range("O3").Value = [double]"0.009994987279658562"
$ws.Range("P3").Value = [double]"0.009994987279658561"
$ws.Range("Q3").Value = [double]"2.411181444901999"
$ws.Range("R3").Value = [double]"21.700633004118"
$ws.Range("S3").Value = [double]"0.0003281670720061659"
$ws.Range("T3").Value = [double]"0.0003281670720061658"
$ws.Range("G4").Value = [double]"29.663058"
$ws.Range("H4").Value = [double]"88.98917399999999"
$ws.Range("I4").Value = [double]"0.03283316554829836"
$ws.Range("J4").Value = [double]"0.03283316554829836"
$ws.Range("M4").Value = [double]"7.912604999999999"
$ws.Range("N4").Value = [double]"23.737815"
$ws.Range("O4").Value = [double]"0.9729438112167713"
$ws.Range("P4").Value = [double]"0.9729438112167712"
$ws.Range("Q4").Value = [double]"234.71206104609"
$ws.Range("R4").Value = [double]"2112.408549414809"
$ws.Range("S4").Value = [double]"0.0319448252228726"
$ws.Range("T4").Value = [double]"0.03194482522287259"
$ws.Range("G5").Value = [double]"29.663058"
$ws.Range("H5").Value = [double]"88.98917399999999"
$ws.Range("I5").Value = [double]"0.03283316554829836"
$ws.Range("J5").Value = [double]"0.03283316554829836"
$ws.Range("K5").Value = [double]"3"
$ws.Range("L5").Value = [double]"1"
$ws.Range("M5").Value = [double]"0.126567"
$ws.Range("N5").Value = [double]"0.379701"
$ws.Range("O5").Value = [double]"0.015562836683276"
$ws.Range("P5").Value = [double]"0.015562836683276"
$ws.Range("Q5").Value = [double]"3.754364261885999"
$ws.Range("R5").Value = [double]"33.78927835697399"
$ws.Range("S5").Value = [double]"0.0005109771932231315"
$ws.Range("T5").Value = [double]"0.0005109771932231313"
$ws.Range("I6").Value = [double]"0.593748363803188"
$ws.Range("J6").Value = [double]"0.5937483638031879"
$ws.Range("O6").Value = [double]"0.001498364820294181"
$ws.Range("P6").Value = [double]"0.001498364820294181"
$ws.Range("S6").Value = [double]"0.0008896516604299276"
$ws.Range("T6").Value = [double]"0.0008896516604299274"
$ws.Range("I7").Value = [double]"0.593748363803188"
$ws.Range("J7").Value = [double]"0.5937483638031879"
$ws.Range("M7").Value = [double]"0.08128566666666666"
$ws.Range("N7").Value = [double]"0.243857"
$ws.Range("O7").Value = [double]"0.009994987279658562"
$ws.Range("P7").Value = [double]"0.009994987279658561"
$ws.Range("Q7").Value = [double]"43.60332041810589"
$ws.Range("R7").Value = [double]"392.429883762953"
$ws.Range("S7").Value = [double]"0.005934507343530948"
$ws.Range("T7").Value = [double]"0.005934507343530946"
$ws.Range("I8").Value = [double]"0.593748363803188"
$ws.Range("J8").Value = [double]"0.5937483638031879"
$ws.Range("M8").Value = [double]"7.912604999999999"
$ws.Range("N8").Value = [double]"23.737815"
$ws.Range("O8").Value = [double]"0.9729438112167713"
$ws.Range("P8").Value = [double]"0.9729438112167712"
$ws.Range("Q8").Value = [double]"4244.485716919015"
$ws.Range("R8").Value = [double]"38200.37145227114"
$ws.Range("S8").Value = [double]"0.5776837959823957"
$ws.Range("T8").Value = [double]"0.5776837959823956"
$ws.Range("I9").Value = [double]"0.593748363803188"
$ws.Range("J9").Value = [double]"0.5937483638031879"
$ws.Range("K9").Value = [double]"3"
$ws.Range("L9").Value = [double]"1"
$ws.Range("M9").Value = [double]"0.126567"
$ws.Range("N9").Value = [double]"0.379701"
$ws.Range("O9").Value = [double]"0.015562836683276"
$ws.Range("P9").Value = [double]"0.015562836683276"
$ws.Range("Q9").Value = [double]"67.893168398181"
$ws.Range("R9").Value = [double]"611.038515583629"
$ws.Range("S9").Value = [double]"0.009240408816831358"
$ws.Range("T9").Value = [double]"0.009240408816831355"
$ws.Range("G10").Value = [double]"54.816723"
$ws.Range("H10").Value = [double]"164.450169"
$ws.Range("I10").Value = [double]"0.06067501675229219"
$ws.Range("J10").Value = [double]"0.06067501675229218"
$ws.Range("O10").Value = [double]"0.001498364820294181"
$ws.Range("P10").Value = [double]"0.001498364820294181"
$ws.Range("Q10").Value = [double]"0.667978314237"
$ws.Range("R10").Value = [double]"6.011804828133"
$ws.Range("S10").Value = [double]"9.091331057239468E-05"
$ws.Range("T10").Value = [double]"9.091331057239468E-05"
$ws.Range("G11").Value = [double]"54.816723"
$ws.Range("H11").Value = [double]"164.450169"
$ws.Range("I11").Value = [double]"0.06067501675229219"
$ws.Range("J11").Value = [double]"0.06067501675229218"
$ws.Range("M11").Value = [double]"0.08128566666666666"
$ws.Range("N11").Value = [double]"0.243857"
$ws.Range("O11").Value = [double]"0.009994987279658562"
$ws.Range("P11").Value = [double]"0.009994987279658561"
$ws.Range("Q11").Value = [double]"4.455813873537"
$ws.Range("R11").Value = [double]"40.102324861833"
$ws.Range("S11").Value = [double]"0.0006064460206322306"
$ws.Range("T11").Value = [double]"0.0006064460206322304"
$ws.Range("G12").Value = [double]"54.816723"
$ws.Range("H12").Value = [double]"164.450169"
$ws.Range("I12").Value = [double]"0.06067501675229219"
$ws.Range("J12").Value = [double]"0.06067501675229218"
$ws.Range("M12").Value = [double]"7.912604999999999"
$ws.Range("N12").Value = [double]"23.737815"
$ws.Range("O12").Value = [double]"0.9729438112167713"
$ws.Range("P12").Value = [double]"0.9729438112167712"
$ws.Range("Q12").Value = [double]"433.743076493415"
$ws.Range("R12").Value = [double]"3903.687688440735"
$ws.Range("S12").Value = [double]"0.0590333820446166"
$ws.Range("T12").Value = [double]"0.05903338204461659"
$ws.Range("G13").Value = [double]"54.816723"
$ws.Range("H13").Value = [double]"164.450169"
$ws.Range("I13").Value = [double]"0.06067501675229219"
$ws.Range("J13").Value = [double]"0.06067501675229218"
$ws.Range("K13").Value = [double]"3"
$ws.Range("L13").Value = [double]"1"
$ws.Range("M13").Value = [double]"0.126567"
$ws.Range("N13").Value = [double]"0.379701"
$ws.Range("O13").Value = [double]"0.015562836683276"
$ws.Range("P13").Value = [double]"0.015562836683276"
$ws.Range("Q13").Value = [double]"6.937988179941"
$ws.Range("R13").Value = [double]"62.441893619469"
$ws.Range("S13").Value = [double]"0.0009442753764709586"
$ws.Range("T13").Value = [double]"0.0009442753764709584"
$ws.Range("G14").Value = [double]"33.372838"
$ws.Range("H14").Value = [double]"100.118514"
$ws.Range("I14").Value = [double]"0.03693941180543633"
$ws.Range("J14").Value = [double]"0.03693941180543633"
$ws.Range("O14").Value = [double]"0.001498364820294181"
$ws.Range("P14").Value = [double]"0.001498364820294181"
$ws.Range("Q14").Value = [double]"0.4066702795886666"
$ws.Range("R14").Value = [double]"3.660032516298"
$ws.Range("S14").Value = [double]"5.534871513162534E-05"
$ws.Range("T14").Value = [double]"5.534871513162534E-05"
$ws.Range("G15").Value = [double]"33.372838"
$ws.Range("H15").Value = [double]"100.118514"
$ws.Range("I15").Value = [double]"0.03693941180543633"
$ws.Range("J15").Value = [double]"0.03693941180543633"
$ws.Range("M15").Value = [double]"0.08128566666666666"
$ws.Range("N15").Value = [double]"0.243857"
$ws.Range("O15").Value = [double]"0.009994987279658562"
$ws.Range("P15").Value = [double]"0.009994987279658561"
$ws.Range("Q15").Value = [double]"2.712733385388666"
$ws.Range("R15").Value = [double]"24.414600468498"
$ws.Range("S15").Value = [double]"0.0003692089511134054"
$ws.Range("T15").Value = [double]"0.0003692089511134054"
$ws.Range("G16").Value = [double]"33.372838"
$ws.Range("H16").Value = [double]"100.118514"
$ws.Range("I16").Value = [double]"0.03693941180543633"
$ws.Range("J16").Value = [double]"0.03693941180543633"
$ws.Range("M16").Value = [double]"7.912604999999999"
$ws.Range("N16").Value = [double]"23.737815"
$ws.Range("O16").Value = [double]"0.9729438112167713"
$ws.Range("P16").Value = [double]"0.9729438112167712"
$ws.Range("Q16").Value = [double]"264.06608482299"
$ws.Range("R16").Value = [double]"2376.59476340691"
$ws.Range("S16").Value = [double]"0.03593997210608701"
$ws.Range("T16").Value = [double]"0.03593997210608701"
$ws.Range("G17").Value = [double]"33.372838"
$ws.Range("H17").Value = [double]"100.118514"
$ws.Range("I17").Value = [double]"0.03693941180543633"
$ws.Range("J17").Value = [double]"0.03693941180543633"
$ws.Range("K17").Value = [double]"3"
$ws.Range("L17").Value = [double]"1"
$ws.Range("M17").Value = [double]"0.126567"
$ws.Range("N17").Value = [double]"0.379701"
$ws.Range("O17").Value = [double]"0.015562836683276"
$ws.Range("P17").Value = [double]"0.015562836683276"
$ws.Range("Q17").Value = [double]"4.223899987146"
$ws.Range("R17").Value = [double]"38.015099884314"
$ws.Range("S17").Value = [double]"0.000574882033104283"
$ws.Range("T17").Value = [double]"0.000574882033104283"
$ws.Range("G18").Value = [double]"205.087789"
$ws.Range("H18").Value = [double]"615.263367"
$ws.Range("I18").Value = [double]"0.2270056353654261"
$ws.Range("J18").Value = [double]"0.227005635365426"
$ws.Range("O18").Value = [double]"0.001498364820294181"
$ws.Range("P18").Value = [double]"0.001498364820294181"
$ws.Range("Q18").Value = [double]"2.499131434157666"
$ws.Range("R18").Value = [double]"22.492182907419"
$ws.Range("S18").Value = [double]"0.000340137258040083"
$ws.Range("T18").Value = [double]"0.0003401372580400829"
$ws.Range("G19").Value = [double]"205.087789"
$ws.Range("H19").Value = [double]"615.263367"
$ws.Range("I19").Value = [double]"0.2270056353654261"
$ws.Range("J19").Value = [double]"0.227005635365426"
$ws.Range("M19").Value = [double]"0.08128566666666666"
$ws.Range("N19").Value = [double]"0.243857"
$ws.Range("O19").Value = [double]"0.009994987279658562"
$ws.Range("P19").Value = [double]"0.009994987279658561"
$ws.Range("Q19").Value = [double]"16.67069765405767"
$ws.Range("R19").Value = [double]"150.036278886519"
$ws.Range("S19").Value = [double]"0.002268918437888243"
$ws.Range("T19").Value = [double]"0.002268918437888243"
$ws.Range("G20").Value = [double]"205.087789"
$ws.Range("H20").Value = [double]"615.263367"
$ws.Range("I20").Value = [double]"0.2270056353654261"
$ws.Range("J20").Value = [double]"0.227005635365426"
$ws.Range("M20").Value = [double]"7.912604999999999"
$ws.Range("N20").Value = [double]"23.737815"
$ws.Range("O20").Value = [double]"0.9729438112167713"
$ws.Range("P20").Value = [double]"0.9729438112167712"
$ws.Range("Q20").Value = [double]"1622.778664680345"
$ws.Range("R20").Value = [double]"14605.0079821231"
$ws.Range("S20").Value = [double]"0.2208637280401224"
$ws.Range("T20").Value = [double]"0.2208637280401223"
$ws.Range("G21").Value = [double]"205.087789"
$ws.Range("H21").Value = [double]"615.263367"
$ws.Range("I21").Value = [double]"0.2270056353654261"
$ws.Range("J21").Value = [double]"0.227005635365426"
$ws.Range("K21").Value = [double]"3"
$ws.Range("L21").Value = [double]"1"
$ws.Range("M21").Value = [double]"0.126567"
$ws.Range("N21").Value = [double]"0.379701"
$ws.Range("O21").Value = [double]"0.015562836683276"
$ws.Range("P21").Value = [double]"0.015562836683276"
$ws.Range("Q21").Value = [double]"25.957346190363"
$ws.Range("R21").Value = [double]"233.616115713267"
$ws.Range("S21").Value = [double]"0.003532851629375429"
$ws.Range("T21").Value = [double]"0.003532851629375428"
$ws.Range("G22").Value = [double]"44.08682333333334"
$ws.Range("H22").Value = [double]"132.26047"
$ws.Range("I22").Value = [double]"0.04879840672535908"
$ws.Range("J22").Value = [double]"0.04879840672535907"
$ws.Range("O22").Value = [double]"0.001498364820294181"
$ws.Range("P22").Value = [double]"0.001498364820294181"
$ws.Range("Q22").Value = [double]"0.5372273335322222"
$ws.Range("R22").Value = [double]"4.835046001789999"
$ws.Range("S22").Value = [double]"7.3117815923685E-05"
$ws.Range("T22").Value = [double]"7.311781592368499E-05"
$ws.Range("G23").Value = [double]"44.08682333333334"
$ws.Range("H23").Value = [double]"132.26047"
$ws.Range("I23").Value = [double]"0.04879840672535908"
$ws.Range("J23").Value = [double]"0.04879840672535907"
$ws.Range("M23").Value = [double]"0.08128566666666666"
$ws.Range("N23").Value = [double]"0.243857"
$ws.Range("O23").Value = [double]"0.009994987279658562"
$ws.Range("P23").Value = [double]"0.009994987279658561"
$ws.Range("Q23").Value = [double]"3.583626825865555"
$ws.Range("R23").Value = [double]"32.25264143279"
$ws.Range("S23").Value = [double]"0.0004877394544875688"
$ws.Range("T23").Value = [double]"0.0004877394544875687"
$ws.Range("G24").Value = [double]"44.08682333333334"
$ws.Range("H24").Value = [double]"132.26047"
$ws.Range("I24").Value = [double]"0.04879840672535908"
$ws.Range("J24").Value = [double]"0.04879840672535907"
$ws.Range("M24").Value = [double]"7.912604999999999"
$ws.Range("N24").Value = [double]"23.737815"
$ws.Range("O24").Value = [double]"0.9729438112167713"
$ws.Range("P24").Value = [double]"0.9729438112167712"
$ws.Range("Q24").Value = [double]"348.84161874145"
$ws.Range("R24").Value = [double]"3139.57456867305"
$ws.Range("S24").Value = [double]"0.04747810782067699"
$ws.Range("T24").Value = [double]"0.04747810782067698"
$ws.Range("G25").Value = [double]"44.08682333333334"
$ws.Range("H25").Value = [double]"132.26047"
$ws.Range("I25").Value = [double]"0.04879840672535908"
$ws.Range("J25").Value = [double]"0.04879840672535907"
$ws.Range("K25").Value = [double]"3"
$ws.Range("L25").Value = [double]"1"
$ws.Range("M25").Value = [double]"0.126567"
$ws.Range("N25").Value = [double]"0.379701"
$ws.Range("O25").Value = [double]"0.015562836683276"
$ws.Range("P25").Value = [double]"0.015562836683276"
$ws.Range("Q25").Value = [double]"5.579936968829999"
$ws.Range("R25").Value = [double]"50.21943271946999"
$ws.Range("S25").Value = [double]"0.0007594416342708406"
$ws.Range("T25").Value = [double]"0.0007594416342708404"

Write-Host "Updated 278 cells"
